$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the header labels in row 2 that previously held placeholder
# "unnamed" text (pandas auto-generated) so that they read "total",
# matching the corrected data used at the start of the PNAD 2009 analysis.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
